# CR_Geometry_WG_Numbers_table.xlsx
# "added 2 WG numbers for the CR"
#
# Appends two new rows (66 and 67) to the "WG NB" sheet, following the
# exact same pattern as all prior rows:
#   A = new WG document number (N92xx)
#   B = author name (already in sharedStrings as index 0)
#   C, D = left blank (just carry the row's border/fill formatting)
#   E = document title / description
#   F = date (already in sharedStrings as index 2 - stored as TEXT, not a
#       real date, matching every other row in the sheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use the last existing data row (65) as the formatting template for the
# two new rows - this reproduces the same per-column styles (s="1", s="1",
# s="5", s="6", s="2", s="1") used throughout the table.
$ws.Range("A65:F65").Copy()
$ws.Range("A66:F67").PasteSpecial(-4122)   # xlPasteFormats

# Fill column A (new WG numbers) first so the new shared-string entries
# are appended in the same order as the source edit.
$ws.Range("A66").Value = "N9213"
$ws.Range("A67").Value = "N9214"

# Column E (the two new publication-set descriptions) next.
$ws.Range("E66").Value = "CR_Geometry wg.number.publication_set"
$ws.Range("E67").Value = "CR_Geometry wg.number.publication_set_comments"

# Column F must stay a plain text date string ("2016-04-15"), not get
# auto-converted to a numeric Excel date - force text format first.
$ws.Range("F66").NumberFormat = "@"
$ws.Range("F67").NumberFormat = "@"

# Columns B and F reuse existing shared strings ("Kevin Le Tutour" and the
# "2016-04-15" date text already used by every other row).
$ws.Range("B66").Value = "Kevin Le Tutour"
$ws.Range("F66").Value = "2016-04-15"
$ws.Range("B67").Value = "Kevin Le Tutour"
$ws.Range("F67").Value = "2016-04-15"

# Re-apply the template formatting to column F so its style index matches
# the rest of the table exactly (setting the value can otherwise leave it
# on the auto-created "text number format" style instead of the shared one).
$ws.Range("F65").Copy()
$ws.Range("F66:F67").PasteSpecial(-4122)

# Scroll/select to mirror the author's on-screen view when the rows were
# added.
$ws.Range("E52").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 41
$win.ScrollColumn = 1
